# Applies the "Resolved Fraud Issue and Added re-implemented CSP Library" edit
# to the Settings sheet of the ESWNikeDailyRefunds Config workbook:
#  - Adds Description (column C) text for several existing settings rows.
#  - Adds a brand new setting row (ConsoleToDateDelay) with Name/Value/Description.
#  - Updates the active selection on the Settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New Description values (column C) for existing rows.
$ws.Range("C6").Value = "Location of Nike Daily Refund Master File"
$ws.Range("C6").WrapText = $true

$ws.Range("C7").Value = "The format of Nike Refund Master File"
$ws.Range("C7").WrapText = $true

$ws.Range("C8").Value = "The location of the UAT Transactions Audit File"
$ws.Range("C8").WrapText = $true

$ws.Range("C9").Value = "The location of the Daily Fraud File"
$ws.Range("C9").WrapText = $true

$ws.Range("C11").Value = "Root Url to CSP"
$ws.Range("C11").WrapText = $true

$ws.Range("C12").Value = "Relative url path to Login Page"
$ws.Range("C12").WrapText = $true

$ws.Range("C13").Value = "Relative url path to search financial transactions page"
$ws.Range("C13").WrapText = $true

$ws.Range("C24").Value = "The phrase to search in the master file for potential fraud transactions"
$ws.Range("C24").WrapText = $true

$ws.Range("C25").Value = "The phrase to search in the master file for potential do not process transactions"
$ws.Range("C25").WrapText = $true

# New row 29: ConsoleToDateDelay setting.
$ws.Range("A29").Value = "ConsoleToDateDelay"
$ws.Range("B29").Value = 500
$ws.Range("C29").Value = "The amount in milliseconds to delay between opening the browser console and entering the jquery statement to update the to and from dates."
$ws.Range("C29").WrapText = $true

# Update the selection shown when the sheet is reopened.
$ws.Activate() | Out-Null
$ws.Range("B21").Select() | Out-Null
